$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Row 2 (Beta) values
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 9.007479029934334
$ws.Range("G2").Value = 8.808168683278657
$ws.Range("H2").Value = 9.210392217009293
$ws.Range("I2").Value = 0.009747342287103821
$ws.Range("J2").Value = 0.008760115367642282
$ws.Range("K2").Value = 0.0109364128083565
$ws.Range("L2").Value = 0.004702444070009214
$ws.Range("M2").Value = 0.004464766308060805
$ws.Range("N2").Value = 0.004973624662091128

# Update existing Row 3 (Gamma) values
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.2773225874354903
$ws.Range("G3").Value = 0.2765055622009212
$ws.Range("H3").Value = 0.278168124700982
$ws.Range("I3").Value = 0.2494149438776253
$ws.Range("J3").Value = 0.2486451159744145
$ws.Range("K3").Value = 0.2502098411367306
$ws.Range("L3").Value = 0.2747273246432781
$ws.Range("M3").Value = 0.273913793545053
$ws.Range("N3").Value = 0.2755698775706886

# Add new Row 4 (Beta + Gamma)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 9.284801617369823
$ws.Range("G4").Value = 9.084674245479578
$ws.Range("H4").Value = 9.488560341710276
$ws.Range("I4").Value = 0.2591622861647291
$ws.Range("J4").Value = 0.2574052313420567
$ws.Range("K4").Value = 0.2611462539450872
$ws.Range("L4").Value = 0.2794297687132873
$ws.Range("M4").Value = 0.2783785598531138
$ws.Range("N4").Value = 0.2805435022327798

# Apply the same style as A2/A3 to A4
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
